$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "Definitions" table (rows 20-24): header row + Micro/Small/Medium/Large
# ---------------------------------------------------------------------------

# Row 20: header row (bold "title" style, like the other table headers)
$ws.Range("B20:D20").Style = "title"
$ws.Range("B20").Value = "Number of employees"
$ws.Range("C20").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D20").Value = "Turnover (local currency, unless noted otherwise)"

# Row 21: Micro
$ws.Range("A21").Value = "Micro"
$ws.Range("B21").Value = "1-20"

# Row 22: Small
$ws.Range("A22").Value = "Small"
$ws.Range("D22").Value = "< limits set at the Global Unique Contribution (CGU) by the General Tax Code"

# Row 23: Medium (A23 previously held the "title"-styled "Min PME - ADEPME"
# text, so its style must be reset back to Normal along with the new value)
$ws.Range("A23").Value = "Medium"
$ws.Range("A23").Style = "Normal"
$ws.Range("B23").Value = "21-250"
$ws.Range("D23").Value = "<= F CFA 5 Billionlion"

# Row 24: Large (A24 previously held the "source"-styled long citation text,
# so its style must be reset back to Normal along with the new value)
$ws.Range("A24").Value = "Large"
$ws.Range("A24").Style = "Normal"
$ws.Range("B24").Value = ">250"
$ws.Range("D24").Value = "> F CFA 5 Billionlion"

# ---------------------------------------------------------------------------
# Source lines for the new table, moved down to rows 29-30
# ---------------------------------------------------------------------------

$ws.Range("A29").Value = "Min PME - ADEPME"
$ws.Range("A29").Style = "title"

$ws.Range("A30").Value = "Ministere des mines, de l'Industrie, de l'Agro-industrie et des PME, Direction des Petites et Moyennes Entreprises (Min PME - ADEPME), ""LETTRE  DE  POLITIQUE SECTORIELLE  DES  PME"", 2010, p. 9. Available at http://www.senegal-entreprises.net/3-download/lettre-politique-sectorielle-10-2010.pdf"
$ws.Range("A30").Style = "source"
